# Actualización automática 2025-11-10 16:30:09
# Updates the "CUMPLIMIENTO MENSUAL" sheet: refresh PRESUPUESTO / VENTA /
# POR CUMPLIR / CUMPLIMIENTO figures and narrow column D slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column D ("VENTA") gets one unit narrower (raw width 14 -> 13).
# Excel's ColumnWidth property is expressed in characters, not the raw
# OOXML "width" units, so 12.14 is the character width that round-trips
# to width="13" in the saved XML.
$ws.Columns.Item(4).ColumnWidth = 12.14

# Row 2 - 240X120 PORCELANATO
$ws.Range("C2").Value = 129.6
$ws.Range("E2").Value = 129.6

# Row 3 - 240X80 PORCELANATO
$ws.Range("C3").Value = 1867.69
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1867.69
$ws.Range("F3").Value = 0

# Row 4 - FREGADEROS DE COCINA
$ws.Range("C4").Value = 1987.7
$ws.Range("D4").Value = 2172.1
$ws.Range("E4").Value = -184.3999999999999
$ws.Range("F4").Value = 1.092770538813704

# Row 5 - GRIFERIAS
$ws.Range("C5").Value = 86.41
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 86.41
$ws.Range("F5").Value = 0

# Row 6 - INODOROS
$ws.Range("C6").Value = 1815
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1815
$ws.Range("F6").Value = 0

# Row 7 - LAVABOS
$ws.Range("C7").Value = 383.4
$ws.Range("E7").Value = 383.4

# Row 8 - NO RESURTIBLES
$ws.Range("C8").Value = 415
$ws.Range("D8").Value = 29.49
$ws.Range("E8").Value = 385.51
$ws.Range("F8").Value = 0.07106024096385542

# Row 10 - PANELES DECORATIVOS
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 388.107983534392
$ws.Range("F10").Value = 0

# Row 11 - PIEDRA SINTERIZADA
$ws.Range("C11").Value = 1440.92
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 1440.92
$ws.Range("F11").Value = 0

# Row 12 - PORCELANATO
$ws.Range("C12").Value = 48041
$ws.Range("D12").Value = 1602.63
$ws.Range("E12").Value = 46438.37
$ws.Range("F12").Value = 0.03335963031577194

# Row 14 - TOTAL
$ws.Range("C14").Value = 57887.35196497848
$ws.Range("D14").Value = 3804.22
$ws.Range("E14").Value = 54083.13196497848
$ws.Range("F14").Value = 0.06571763728804751
